$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old data validation + old helper/lookup cell content -------------
$ws.Range("F3:K10").Validation.Delete()
$ws.Range("F3").ClearContents()
$ws.Range("F4:G8").ClearContents()
$ws.Range("K4:K7").ClearContents()
$ws.Range("H16:H18").ClearContents()

# --- New header row (row 2): Iteration / Preprocessing / existing headers ----
$ws.Range("A2").Value = "Iteration"
$ws.Range("B2").Value = "Preprocessing"
# C2, D2, E2 already hold the right text - just need restyling below

# Header formatting: reuse the existing bold+bottom-border style (currently on
# C2:E2) by copying it across the header, then add center + wrap-text on top.
$ws.Range("C2").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$ws.Range("A2:E2").HorizontalAlignment = -4108
$ws.Range("A2:E2").WrapText = $true
$ws.Application.CutCopyMode = $false

# F2: bold + centered + wrapped but no border, kept blank
$ws.Range("C2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Borders.Item(9).LineStyle = -4142
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").WrapText = $true
$ws.Application.CutCopyMode = $false

# Row 2 is taller to fit the wrapped header text
$ws.Rows.Item(2).RowHeight = 45

# --- New data columns D (# of Features) and E (Log-Likelihood) ---------------
# (Column A's index cells 3-8 already carry the centered style.)
$ws.Range("D3").Value = 49
$ws.Range("E3").Value = 603.81
$ws.Range("D4").Value = 47
$ws.Range("E4").Value = 466.13
$ws.Range("D5").Value = 45
$ws.Range("E5").Value = 516.28
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 466.56
$ws.Range("D7").Value = 44
$ws.Range("E7").Value = 498.82
$ws.Range("D8").Value = 41
$ws.Range("E8").Value = 396.41
$ws.Range("D3:E8").HorizontalAlignment = -4108

# --- Move the footnote list from column H to column G ------------------------
$ws.Range("G16").Value = "Removed podcast stop words"
$ws.Range("G17").Value = "removed short words"
$ws.Range("G18").Value = "iTunes categories weren't well set up. Usubs found junk and then better categorized existing categories. For example …"

# --- Column widths -------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.7109375
$ws.Columns.Item(4).ColumnWidth = 21.5703125
$ws.Columns.Item(5).ColumnWidth = 23.7109375

# --- Selection lands on B2 like the authored workbook -------------------------
$ws.Range("B2").Select()
